# Convert the job-array "years" in column B (on both sheets) into small
# sequential numbers (year - 2000) so the concatenated job numbers stay
# under 1000, per the commit message "changing so that the job nums were
# less than 1000".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("option 1")
$ws2 = $wb.Worksheets.Item("option 2")

$years = @(2001,2002,2003,2004,2005,2006,2007,2008,2009,2010,2011,2015,2016,2017,2018,2019,2020,2021,2022)
$rows  = @(4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $newVal = $years[$i] - 2000
    $ws1.Cells.Item($r, 2).Value = $newVal
    $ws2.Cells.Item($r, 2).Value = $newVal
}

# Restore the view/selection state recorded in the saved workbook: sheet
# "option 2" keeps a selection on B27 (no longer the active tab), and
# sheet "option 1" becomes the active tab with the selection on S16.
$ws2.Range("B27").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("S16").Select() | Out-Null
